# "new groups can now be created"
# - Flip the "Actief" (active) flag: Totaal becomes inactive, test becomes active
# - Make room for new groups by extending the table with 5 blank rows
# - Leave the selection on the first new blank row (A6), ready for data entry

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle the Actief column values
$ws.Range("I2").Value = $false
$ws.Range("I3").Value = $true

# Extend the used range down to row 8 so new groups have room to be added
$ws.Range("A4:I8").Style = "Normal"

# Position the selection on the first blank row, ready for a new group
$ws.Range("A6").Select()
